$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

# Column C ("Förändrad") holds a date serial number that was bumped by one day
# (46060 -> 46061) for every data row (rows 2 through 147).
$ws.Range("C2:C147").Value = 46061
